$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose "Sending cluster" is ECs (original rows 2-4).
# Deleting row index 2 three times in a row shifts the remaining rows
# (formerly 5-7, Sending cluster = MuSCs) up into rows 2-4.
$ws.Rows.Item(2).EntireRow.Delete()
$ws.Rows.Item(2).EntireRow.Delete()
$ws.Rows.Item(2).EntireRow.Delete()

# Refresh the surviving MuSCs-sourced rows (now rows 2-4) with the newly
# recomputed TPM-based values.
$rowData = @{
  2 = @{
    Labels = @("MuSCs", "Il1rapl1", "Ptprf", "ECs")
    Nums   = @(3, 1, 0.257516, 0.772548, 1, 1, 3, 1, 0.05031533333333333, 0.150946, 0.005485022167780355, 0.005485022167780356, 0.01295700337866667, 0.116613030408, 0.005485022167780355, 0.005485022167780356)
  }
  3 = @{
    Labels = @("MuSCs", "Il1rapl1", "Ptprf", "FAPs")
    Nums   = @(3, 1, 0.257516, 0.772548, 1, 1, 3, 1, 3.467027333333334, 10.401082, 0.377950825718477, 0.377950825718477, 0.8928150107706668, 8.035335096936, 0.377950825718477, 0.377950825718477)
  }
  4 = @{
    Labels = @("MuSCs", "Il1rapl1", "Ptprf", "MuSCs")
    Nums   = @(3, 1, 0.257516, 0.772548, 1, 1, 3, 1, 5.655880666666666, 16.967642, 0.6165641521137426, 0.6165641521137426, 1.456479765757333, 13.108317891816, 0.6165641521137426, 0.6165641521137426)
  }
}

foreach ($r in @(2, 3, 4)) {
  $labels = $rowData[$r].Labels
  for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($r, 1 + $i).Value = $labels[$i]
  }
  $nums = $rowData[$r].Nums
  for ($i = 0; $i -lt $nums.Length; $i++) {
    $ws.Cells.Item($r, 5 + $i).Value = $nums[$i]
  }
}
